# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Tomate" (Vega Monumental Concepción)
# just above the existing row 422, shifting the remaining historical rows
# down by two (422-444 -> 424-446), and populate the two new rows with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 422, pushing old rows 422:444 down to 424:446
$ws.Rows("422:423").Insert()

# --- New row 422 ---
$ws.Cells.Item(422,1).Value = 11
$ws.Cells.Item(422,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(422,3).Value = "Bíobío"
$ws.Cells.Item(422,4).Value = 44747
$ws.Cells.Item(422,5).Value = 8
$ws.Cells.Item(422,6).Value = 100112020
$ws.Cells.Item(422,7).Value = "Tomate"
$ws.Cells.Item(422,8).Value = "Larga vida"
$ws.Cells.Item(422,9).Value = "Extra"
$ws.Cells.Item(422,10).Value = 350
$ws.Cells.Item(422,11).Value = 14000
$ws.Cells.Item(422,12).Value = 15000
$ws.Cells.Item(422,13).Value = 14429
$ws.Cells.Item(422,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(422,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(422,16).Value = 802
$ws.Cells.Item(422,17).Value = 18
$ws.Cells.Item(422,18).Value = "Hortaliza"

# --- New row 423 ---
$ws.Cells.Item(423,1).Value = 11
$ws.Cells.Item(423,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(423,3).Value = "Bíobío"
$ws.Cells.Item(423,4).Value = 44747
$ws.Cells.Item(423,5).Value = 8
$ws.Cells.Item(423,6).Value = 100112020
$ws.Cells.Item(423,7).Value = "Tomate"
$ws.Cells.Item(423,8).Value = "Larga vida"
$ws.Cells.Item(423,9).Value = "Primera"
$ws.Cells.Item(423,10).Value = 550
$ws.Cells.Item(423,11).Value = 11000
$ws.Cells.Item(423,12).Value = 12000
$ws.Cells.Item(423,13).Value = 11545
$ws.Cells.Item(423,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(423,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(423,16).Value = 641
$ws.Cells.Item(423,17).Value = 18
$ws.Cells.Item(423,18).Value = "Hortaliza"
